# Apply required changes to Sports output file (Golf sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $typeCell = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $oldType = $typeCell.Value2
    $oldValue = $valueCell.Value2

    if ($oldType -eq $null -or $oldValue -eq $null) {
        continue
    }

    # Determine prefix (club / uil) from the "type" column, e.g. "club-sports" -> "club"
    $prefix = $oldType.Split("-")[0]

    # Determine gender suffix from the "value" column, e.g. "Golf-Boys" -> "boys"
    $parts = $oldValue.Split("-")
    if ($parts.Length -ge 2) {
        $gender = $parts[1].ToLower()
        $newType = "sports_" + $prefix + "_" + $gender
        $newValue = $parts[0]

        $typeCell.Value = $newType
        $valueCell.Value = $newValue
    }
}
